$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the extent of the existing table before we start overwriting it,
# so we know how far down any now-obsolete rows reach.
$oldLastRow = $ws.UsedRange.Rows.Count

$data = @(
    @('574 Plummer Street Building B1 & B2 Port Melbourne', 5, 'old'),
    @('A1 Bakery Brunswick', 5, 'new'),
    @('A1 Bakery Brunswick', 6, 'old'),
    @('Acquire BPO Southbank', 12, 'new'),
    @('Acquire BPO Southbank', 15, 'old'),
    @('Al Haj Halal Meats Glenroy', 52, 'old'),
    @('Al Haj Halal Meats Glenroy', 53, 'new'),
    @('Al-Taqwa College Truganina', 12, 'new'),
    @('Al-Taqwa College Truganina', 14, 'old'),
    @('Best&Less Fountain Gate Narre Warren', 6, 'old'),
    @('Best&Less Fountain Gate Narre Warren', 7, 'new'),
    @('Can Panel Cambellfield', 6, 'old'),
    @('Can Panel Campbellfield', 5, 'new'),
    @('Caroline Springs Police Station', 7, 'old'),
    @('Caroline Springs Police Station', 8, 'new'),
    @('Cedars Medical Clinic Coburg', 44, 'new'),
    @('Cedars Medical Clinic Coburg', 49, 'old'),
    @('Central Plaster 236 & 259 Jasper Road McKinnon', 6, 'old'),
    @('Chemist Warehouse Fillo Drive Somerton', 5, 'new'),
    @('City of Moreland Community', 6, 'new'),
    @('City of Moreland Community', 7, 'old'),
    @('City of Wyndham Community', 5, 'old'),
    @('Coles Aurora Village Epping', 5, 'new'),
    @('Coles Campbellfield Plaza Campbellfield', 8, 'new'),
    @('Coles Campbellfield Plaza Campbellfield', 9, 'old'),
    @('Coles Coburg North Village', 14, 'old'),
    @('Coles Coburg North Village', 15, 'new'),
    @('Coles Roxburgh Village Roxburgh Park', 11, 'new'),
    @('Coles Roxburgh Village Roxburgh Park', 15, 'old'),
    @('Community Kids Meadow Heights', 19, 'new'),
    @('Community Kids Meadow Heights', 21, 'old'),
    @('Construction Site Olea Apartment Caulfield North', 5, 'old'),
    @('Construction Site Olea Apartment Caulfield North', 8, 'new'),
    @('Costco Wholesale Epping', 21, 'new'),
    @('Costco Wholesale Epping', 22, 'old'),
    @('Crossroads Logistics Sunshine North', 5, 'old'),
    @('Croydon Orthodontics', 5, 'new'),
    @('Croydon Orthodontics', 6, 'old'),
    @('Direct Freight Express Campbellield', 6, 'new'),
    @('Fitzroy Community School Fitzroy North', 55, 'new'),
    @('Fitzroy Community School Fitzroy North', 57, 'old'),
    @('Gloria Jeans Coffees Broadmeadows Central', 5, 'old'),
    @('Hamilton Marino 236 Jasper Road McKinnon', 6, 'new'),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 5, 'old'),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 6, 'new'),
    @('Ilim Leaning Sanctuary Glenroy', 15, 'new'),
    @('Ilim Learning Sanctuary Glenroy', 18, 'old'),
    @('Islamic College of Melbourne Tarneit', 8, 'new'),
    @('Islamic College of Melbourne Tarneit', 9, 'old'),
    @('KFC Fawkner', 11, 'new'),
    @('KFC Fawkner', 14, 'old'),
    @('Kasr Sweets Coolaroo', 9, 'new'),
    @('Kasr Sweets Coolaroo', 11, 'old'),
    @('Learning Nest Early Learning Centre Meadow Heights', 9, 'new'),
    @('Learning Nest Early Learning Centre Meadow Heights', 10, 'old'),
    @('Level Crossing Removal Project Lilydale Construction Site John Street', 8, 'new'),
    @('Level Crossing Removal Project Lilydale Construction Site John Street', 9, 'old'),
    @('Lineage Logistics Laverton North', 7, 'new'),
    @('Linfox Somerton National Distribution Centre Somerton', 6, 'new'),
    @('Linfox Somerton National Distribution Centre Somerton', 7, 'old'),
    @('Malvern Health and Fitness Clinic Malvern', 5, 'new'),
    @('Malvern Health and Fitness Clinic Malvern', 6, 'old'),
    @('McDonalds Thomastown II', 6, 'new'),
    @('Melton Police Station Melton', 5, 'new'),
    @('North Geelong House Party', 7, 'new'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 48, 'old'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 52, 'new'),
    @('OnQ Plumbing and Excavations Craigieburn', 11, 'old'),
    @('OnQ Plumbing and Excavations Craigieburn', 15, 'new'),
    @('Oporto Coolaroo', 15, 'old'),
    @('Oporto Coolaroo', 17, 'new'),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 28, 'new'),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 40, 'old'),
    @('Ramsay Health Care Warringal Private Hospital Heidelberg', 9, 'new'),
    @('Serco Mill Park', 5, 'new'),
    @('Serco Mill Park', 6, 'old'),
    @('Southern Cross Station Crew Room Tissue Box Docklands', 5, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne', 5, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne', 6, 'new'),
    @('Tek Foods Somerton', 10, 'old'),
    @('Tek Foods Somerton', 11, 'new'),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 6, 'old'),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 7, 'new'),
    @('Total Window Concepts Hoppers Crossing', 6, 'old'),
    @('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 10, 'old'),
    @('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 13, 'new'),
    @('Victorian Civil and Administrative Tribunal Melbourne', 5, 'old'),
    @('Werribee Mercy Hospistal Emergency Department', 7, 'new'),
    @('Western Health Sunshine Hospital Emergency Department', 5, 'old'),
    @('Western Health Sunshine Hospital Emergency Department', 7, 'new'),
    @('Woolworths Greenvale Lakes Roxburgh Park', 10, 'new'),
    @('Woolworths Greenvale Lakes Roxburgh Park', 13, 'old'),
    @('Yara Childcare Centre Truganina', 26, 'old'),
    @('Yarra Childcare Centre Truganina', 15, 'new')
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $item = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $item[0]
    $ws.Cells.Item($r, 2).Value2 = $item[1]
    $ws.Cells.Item($r, 3).Value2 = $item[2]
}

# Remove old trailing rows that no longer exist in the refreshed data set
$lastRow = $rowCount + 2
if ($oldLastRow -ge $lastRow) {
    $clearRange = $ws.Range("A" + $lastRow + ":C" + $oldLastRow)
    $clearRange.ClearContents()
}

Write-Host "Updated" $rowCount "rows"
